# Add the "2020" column (column N) of data to the 6.4.1.2 water-loss
# indicator sheet, mirroring the formatting already used by the
# neighbouring "2019" column (column M), and leave the selection where
# the author left it (M25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the column-M number formatting into column N -----------------
# Row 4 is an intentionally blank separator row between the two data
# blocks and is left untouched (it never gets an "N" cell), so the
# format copy is split in two so it isn't dragged along for the ride.
$ws.Range("M3").Copy()
$ws.Range("N3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("M5:M25").Copy()
$ws.Range("N5:N25").PasteSpecial(-4122)   # xlPasteFormats

# --- Fill in the 2020 values --------------------------------------------
$ws.Range("N3").Value = 2020

$ws.Range("N5").Value = 2198.6999999999998
$ws.Range("N6").Value = 132.69999999999999
$ws.Range("N7").Value = 242.9
$ws.Range("N8").Value = 203.3
$ws.Range("N9").Value = 202.8
$ws.Range("N10").Value = 284.7
$ws.Range("N11").Value = 294.89999999999998
$ws.Range("N12").Value = 802.5
$ws.Range("N13").Value = 28.1
$ws.Range("N14").Value = 6.8

$ws.Range("N16").Value = 27.4
$ws.Range("N17").Value = 17.5
$ws.Range("N18").Value = 24.7
$ws.Range("N19").Value = 31.5
$ws.Range("N20").Value = 30.4
$ws.Range("N21").Value = 24.8
$ws.Range("N22").Value = 30.7
$ws.Range("N23").Value = 30.1
$ws.Range("N24").Value = 21.2
$ws.Range("N25").Value = 11.6

# --- Restore the author's final selection -------------------------------
$ws.Range("M25").Select()
